# Adds Artisan Commands `showCurve`, `showExtraCurve`, `showEvents`, and
# `showBackgroundEvents` to the "Commands" sheet, right before the existing
# "RC Command" section (which - together with everything below it - shifts
# down by four rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# Insert four new blank rows right above the "RC Command" row (row 96).
$ws.Rows.Item(96).Insert()
$ws.Rows.Item(96).Insert()
$ws.Rows.Item(96).Insert()
$ws.Rows.Item(96).Insert()

# Fill in the four new Command/Description pairs.
$ws.Cells.Item(96, 2).Value2 = "showCurve(<name>,<bool>)"
$ws.Cells.Item(96, 3).Value2 = "shows/hides the curve indicated by <name> which is one of { ET, BT, DeltaET, DeltaBT, BackgroundET, BackgroundBT}"

$ws.Cells.Item(97, 2).Value2 = "showExtraCurve(<extra_device>,<curve>,<bool>)"
$ws.Cells.Item(97, 3).Value2 = "shows/hides the <curve> (one of {T1,T2}) of the zero-based <extra_device> number"

$ws.Cells.Item(98, 2).Value2 = "showEvents(<event_type>, <bool>)"
$ws.Cells.Item(98, 3).Value2 = "shows/hides the events of <event_type> in [1,..,5]"

$ws.Cells.Item(99, 2).Value2 = "showBackgroundEvents(<bool>)"
$ws.Cells.Item(99, 3).Value2 = "shows/hides the events of the background profile"

# Restore the "Commands" sheet as the active tab with C97 selected, matching
# the author's final cursor position after the edit.
$ws.Activate()
$ws.Range("C97").Select()
